# Update "想去人数" (want-to-go count) figures across sheets to match
# the freshly generated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 71
$wsExpo.Range("F5").Value = 4676
$wsExpo.Range("F6").Value = 365
$wsExpo.Range("F10").Value = 206

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 20

# Sheet "全部类型" (All types) - combined listing
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 71
$wsAll.Range("F5").Value = 4676
$wsAll.Range("F6").Value = 365
$wsAll.Range("F10").Value = 20
$wsAll.Range("F11").Value = 206
